# Updates odds values in "Sheet1" for the two matches in row 6 and row 17
# (Flashscore odds refresh). Only numeric odds cells change; all labels,
# dates, teams, leagues and styling stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 6: East Bengal vs Punjab (INDIA - ISL) ----
$ws.Range("G6").Value = 2.4
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.6
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 2.2
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 13
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 4
$ws.Range("Q6").Value = 1.75
$ws.Range("R6").Value = 2.05
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("W6").Value = 9.5
$ws.Range("Y6").Value = 10
$ws.Range("Z6").Value = 23
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 41
$ws.Range("AG6").Value = 151
$ws.Range("AH6").Value = 10
$ws.Range("AI6").Value = 15
$ws.Range("AJ6").Value = 10
$ws.Range("AK6").Value = 26
$ws.Range("AL6").Value = 21
$ws.Range("AM6").Value = 26
$ws.Range("AN6").Value = 4.75
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 41
$ws.Range("AR6").Value = 51
$ws.Range("AS6").Value = 126
$ws.Range("AT6").Value = 3
$ws.Range("AU6").Value = 7.5
$ws.Range("AV6").Value = 51
$ws.Range("AX6").Value = 15
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 41
$ws.Range("BA6").Value = 67
$ws.Range("BB6").Value = 151

# ---- Row 17: Shabab Al-Ahli Dubai vs Al Wasl (UAE LEAGUE) ----
$ws.Range("G17").Value = 1.9
$ws.Range("H17").Value = 4.05
$ws.Range("I17").Value = 3.25
$ws.Range("J17").Value = 2.32
$ws.Range("L17").Value = 3.5
$ws.Range("P17").Value = 5.6
$ws.Range("W17").Value = 15
$ws.Range("Z17").Value = 20
$ws.Range("AB17").Value = 15.5
$ws.Range("AQ17").Value = 26
$ws.Range("AY17").Value = 17
$ws.Range("AZ17").Value = 65
